$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Updates")

for ($r = 5; $r -le 96; $r++) {
    $ws.Cells.Item($r, 2).Value = 43932
}

$text156 = @'
• The EUR 37 billion Coronavirus Response Investment Initiative (CRII; proposed on 13 March, in force since 1 April) will accelerate the implementation of cohesion policy (total country allocations for the 2014-2020 period remain unchanged). Member States will not have to pay back to the EU budget unspent pre-financing received from structural funds (about EUR 8 billion). Instead, they will be able to use it as national co-financing for the next tranches of their structural fund allocations (about EUR 29 billion). Investment for fostering the crisis response capacities in public health services and some financing of working capital in SMEs will become eligible.
• The Coronavirus Response Investment Initiative Plus (CRII+), proposed on 2 April, complements the CRII by further enhancing flexibility in the use of cohesion funds. This enhanced flexibility is inter alia provided through transfer possibilities across the three cohesion policy funds (the European Regional Development Fund, European Social Fund and Cohesion Fund), transfers between the different categories of regions (e.g. less vs more developed), flexibility regarding thematic concentration, the possibility for a 100% EU co-financing rate for the accounting year 2020-2021, and simplified procedural steps.
• Maximum flexibility will be deployed in the application of the Stability and Growth Pact by activating for the first time (on 23 March) the general escape clause in case of a severe economic downturn for the euro area or EU as a whole.
• The scope of the EU Solidarity Fund has been extended to cover major public health emergencies and possible advance payments from this Fund have been increased (proposal on 13 March, in force since 1 April). Up to EUR 800 million is available in 2020.
• At the request of Member States, the European Globalisation Adjustment Fund could be mobilised to support dismissed workers and those self-employed (up to EUR 179 million available in 2020)
• 	EUR 1 billion has been provided as a guarantee to the European Investment Fund (part of the EIB group) to support SME financing (announced on 13 March, implemented on 6 April). 
• On 9 April, EU finance ministers decided to establish Pandemic Crisis Support credit lines within the framework of the European Stability Mechanism (ESM). Access granted will be 2% of the respective country’s GDP as of end-2019, as a benchmark (about €240 billion in total). The credit line will be available until the COVID 19 crisis is over. The only requirement to access the credit line is that euro area Member States requesting support would commit to use this credit line to finance direct and indirect healthcare, cure and prevention related costs due to the COVID 19 crisis.
'@
$text157 = @'
• The Commission intends to allow State aid for struggling businesses. In particular, the Commission considers that the impact of the COVID-19 outbreak is causing a serious disturbance to the entire EU economy (Article 107(3)(b) of the Treaty). In this context, the Commission adopted on 19 March a Temporary Framework, in place until end-2020, to enable Member States to use the full flexibility foreseen under State aid rules.
• The European Investment Bank (EIB) Group announced on 16 March a plan (currently being deployed) to mobilize up to EUR 40 billion of financing to alleviate liquidity and working capital constraints for SMEs and mid-caps. This comprises dedicated guarantee schemes to banks based on existing programmes, dedicated liquidity lines to banks, and dedicated asset-backed securities (ABS) purchasing programmes to allow banks to transfer risk on portfolios of SME loans.
• In addition to the above plan, EU finance ministers endorsed on 9 April an EIB proposal to create a EUR 25 billion guarantee fund, which will support up to EUR 200 billion of financing for companies (especially SMEs) throughout the EU. The scheme will be implemented by the EIB Group, in close partnership with national promotional banks and other financial intermediaries. 
• The Commission presented on 2 April a proposal for a new instrument for temporary Support to mitigate Unemployment Risks in an Emergency (SURE). SURE support will take the form of loans granted on favourable terms from the EU to Member States, to help them cover the costs directly related to the creation or extension of national short-time work schemes, and other similar measures for the self-employed, in the context of the current crisis. EUR 100 billion (0.7% of 2019 EU27 GDP) will be available for this instrument (with no pre-allocated national envelopes), backed by EUR 25 billion of guarantees voluntarily committed by Member States to the EU budget. SURE will have a temporary nature: its duration and scope are limited to tackling the consequences of the coronavirus pandemic. An instrument along these lines was endorsed by EU finance ministers on 9 April.
'@
$text158 = @'
• The ECB took several measures to support bank lending and liquidity throughout the euro area:
o 12 March: lowering the interest rate applied in targeted longer-term refinancing operations (TLTRO III) during the period from June 2020 to June 2021 (25 basis points below the average rate applied in the Eurosystem’s main refinancing operations).
o 12 March: conducting additional longer-term refinancing operations (LTROs) temporarily (with an interest rate equal to the average rate on the deposit facility, -0.50% currently).
o 18 March: easing collateral standards to give easier access to ECB liquidity, by adjusting the main risk parameters of the collateral framework.
o Major banks (i.e. those directly supervised by the ECB) will be allowed to operate temporarily below the level of capital defined by the Pillar 2 Guidance (P2G), the capital conservation buffer (CCB) and the liquidity coverage ratio (LCR). The ECB considers that these temporary measures will be enhanced by the appropriate relaxation of the countercyclical capital buffer (CCyB) by the national macroprudential authorities (12 March). In addition, the ECB has introduced supervisory flexibility regarding the treatment of non-performing loans (NPLs), in particular to allow banks to fully benefit from public guarantees and moratoriums. Flexibility will also extend to banks’ implementation of NPL reduction strategies (20 March). Furthermore, the ECB has asked banks not to pay dividends until at least 1 October 2020 (27 March).
• The ECB also took measures to reinforce the asset purchase programme (APP):
o 12 March: adding a temporary envelope of additional net asset purchases of €120 billion until the end of the year.
o 18 March: furthermore, launching a new temporary asset purchase programme of private and public sector securities (Pandemic Emergency Purchase Programme, PEPP) with an overall envelope of EUR 750 billion until the end of 2020. Together, both envelopes amount to 7.3% of euro area GDP. Some self-imposed purchase limits will not apply to the PEPP.
o 18 March: expanding the range of eligible assets under the corporate sector purchase programme (CSPP) to non-financial commercial paper. 

'@

$ws.Range("H19").Value = $text156
$ws.Range("J19").Value = $text157
$ws.Range("K19").Value = $text158

Write-Output "done"
